$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Map of row -> column letter -> new value
$changes = @{
    2  = @{ E = 120 }
    4  = @{ E = 53 }
    5  = @{ E = 177; F = 125; H = 136 }
    6  = @{ E = 57; F = 39; H = 49 }
    8  = @{ E = 11 }
    10 = @{ E = 779; F = 461; H = 556 }
    11 = @{ E = 516; F = 315; H = 380 }
    12 = @{ E = 809; F = 518; H = 604 }
    13 = @{ E = 180 }
    14 = @{ E = 160; F = 93; H = 127 }
    15 = @{ E = 221; F = 108; H = 159 }
    16 = @{ E = 246; F = 148; H = 196 }
    17 = @{ E = 138 }
    20 = @{ E = 109 }
    21 = @{ E = 158; F = 97; H = 128 }
    22 = @{ E = 210 }
    23 = @{ E = 249; F = 130; H = 182 }
    24 = @{ E = 308; F = 178; H = 208 }
    25 = @{ E = 365; F = 212; H = 272 }
    26 = @{ E = 239 }
    27 = @{ E = 431; F = 252; H = 334 }
    28 = @{ E = 248; F = 125; H = 177 }
    30 = @{ E = 280; F = 180; H = 233 }
    31 = @{ F = 42; H = 69 }
    32 = @{ E = 230 }
    33 = @{ E = 366; F = 198; H = 289 }
    34 = @{ E = 280; F = 200; H = 238 }
    35 = @{ E = 196; F = 140; H = 167 }
    36 = @{ E = 96; F = 61; H = 71 }
    37 = @{ E = 209; F = 120; H = 156 }
    38 = @{ F = 73; H = 90 }
    39 = @{ E = 216 }
    40 = @{ E = 331; F = 180; H = 260 }
    41 = @{ E = 476; F = 255; H = 347 }
    42 = @{ E = 501 }
    43 = @{ E = 157; F = 96; H = 123 }
    44 = @{ E = 420; F = 237; H = 305 }
    45 = @{ E = 195 }
    46 = @{ E = 420; F = 254; H = 318 }
    47 = @{ E = 597; F = 351; H = 443 }
    48 = @{ E = 299; F = 153; H = 197 }
    49 = @{ E = 362; F = 187; H = 274 }
    50 = @{ E = 301 }
    51 = @{ E = 273; F = 144; H = 218 }
    52 = @{ E = 35 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}

$wb.Save()
